$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.992.66"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").Value = "2.212.42"
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "288.84"
$ws.Range("E5").Value = "  -3.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.93"
$ws.Range("E6").Value = "  +4.56%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +0.93%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.68"
$ws.Range("E10").Value = "  +3.38%  "

$ws.Range("E11").Value = "  -0.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.84"
$ws.Range("E12").Value = "  +3.60%  "

$ws.Range("E13").Value = "  +2.27%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.45"
$ws.Range("E14").Value = "  +2.63%  "

$ws.Range("D15").Value = "2.554.87"
$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.01"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").Value = "2.188.23"
$ws.Range("E17").Value = "  -1.31%  "

$ws.Range("E18").Value = "  +1.39%  "

$ws.Range("D19").Value = "39.930.91"
$ws.Range("E19").Value = "  +0.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.85"
$ws.Range("E20").Value = "  +13.86%  "

$ws.Range("D21").Value = "0.0₃0884"
$ws.Range("E21").Value = "  +0.75%  "

$ws.Range("E22").Value = "  +0.68%  "

$ws.Range("E23").Value = "  +0.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.06"
$ws.Range("E24").Value = "  +0.49%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.46"
$ws.Range("E26").Value = "  +1.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.84"
$ws.Range("E27").Value = "  +0.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.58"
$ws.Range("E28").Value = "  -0.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  +1.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.25"
$ws.Range("E30").Value = "  +0.67%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.73"
$ws.Range("E31").Value = "  +2.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.15"
$ws.Range("E32").Value = "  -0.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.11%  "

$ws.Range("E34").Value = "  +2.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0719"
$ws.Range("E35").Value = "  +2.31%  "

$ws.Range("E36").Value = "  -0.62%  "

$ws.Range("E37").Value = "  +6.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.00"
$ws.Range("E38").Value = "  -1.90%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.112"
$ws.Range("E39").Value = "  +0.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.100"
$ws.Range("E40").Value = "  +2.93%  "

$ws.Range("E41").Value = "  +2.77%  "

$ws.Range("D42").Value = "2.087.03"
$ws.Range("E42").Value = "  +8.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.82"
$ws.Range("E43").Value = "  +4.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.21"
$ws.Range("E44").Value = "  +2.55%  "

$ws.Range("E45").Value = "  +1.01%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.74"
$ws.Range("E46").Value = "  +7.07%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.86"
$ws.Range("E47").Value = "  +6.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.66"
$ws.Range("E48").Value = "  +2.48%  "

$ws.Range("D49").Value = "2.426.20"
$ws.Range("E49").Value = "  -0.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "69.50"
$ws.Range("E50").Value = "  -1.97%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "88.61"
$ws.Range("E51").Value = "  -0.03%  "
